$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.187.94"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3
$ws.Range("D3").Value = "1.786.75"
$ws.Range("E3").Value = "  -0.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").Value = "'226.09"
$ws.Range("E5").Value = "  -0.28%  "

# Row 6
$ws.Range("E6").Value = "  +0.45%  "

# Row 7
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("D8").Value = "'32.03"
$ws.Range("E8").Value = "  -0.58%  "

# Row 9
$ws.Range("E9").Value = "  -0.42%  "

# Row 10
$ws.Range("D10").Value = "'0.0687"
$ws.Range("E10").Value = "  +0.08%  "

# Row 11
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("D12").Value = "2.043.36"
$ws.Range("E12").Value = "  -0.20%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'10.99"
$ws.Range("E13").Value = "  -3.24%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.781.74"
$ws.Range("E14").Value = "  -0.34%  "

# Row 15
$ws.Range("D15").Value = "'0.626"
$ws.Range("E15").Value = "  +0.64%  "

# Row 16
$ws.Range("D16").Value = "34.163.28"
$ws.Range("E16").Value = "  +0.27%  "

# Row 17
$ws.Range("E17").Value = "  +0.34%  "

# Row 18
$ws.Range("D18").Value = "'67.77"
$ws.Range("E18").Value = "  -0.38%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  +2.43%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'245.58"
$ws.Range("E20").Value = "  +0.81%  "

# Row 21
$ws.Range("D21").Value = "'11.02"
$ws.Range("E21").Value = "  +1.05%  "

# Row 22
$ws.Range("E22").Value = "  +0.29%  "

# Row 23
$ws.Range("E23").Value = "  +0.86%  "

# Row 24
$ws.Range("E24").Value = "  -0.50%  "

# Row 26
$ws.Range("E26").Value = "  -0.29%  "

# Row 27
$ws.Range("D27").Value = "'16.28"
$ws.Range("E27").Value = "  +0.15%  "

# Row 28
$ws.Range("D28").Value = "'0.115"
$ws.Range("E28").Value = "  +1.14%  "

# Row 29
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +0.35%  "

# Row 30
$ws.Range("E30").Value = "  -0.23%  "

# Row 31
$ws.Range("E31").Value = "  -0.76%  "

# Row 32
$ws.Range("E32").Value = "  +2.46%  "

# Row 33
$ws.Range("E33").Value = "  +3.28%  "

# Row 34
$ws.Range("E34").Value = "  -2.05%  "

# Row 35
$ws.Range("D35").Value = "1.446.71"
$ws.Range("E35").Value = "  +3.15%  "

# Row 36
$ws.Range("D36").Value = "'2.64"
$ws.Range("E36").Value = "  +12.34%  "

# Row 37
$ws.Range("E37").Value = "  +1.10%  "

# Row 38
$ws.Range("E38").Value = "  +0.94%  "

# Row 39
$ws.Range("E39").Value = "  +0.13%  "

# Row 40
$ws.Range("D40").Value = "'81.86"
$ws.Range("E40").Value = "  +2.01%  "

# Row 41
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.39"
$ws.Range("E41").Value = "  +1.66%  "

# Row 42
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'14.04"
$ws.Range("E42").Value = "  +4.98%  "

# Row 43
$ws.Range("E43").Value = "  +0.81%  "

# Row 44
$ws.Range("D44").Value = "'0.918"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45
$ws.Range("D45").Value = "'0.0516"
$ws.Range("E45").Value = "  +1.34%  "

# Row 46
$ws.Range("D46").Value = "'6.10"
$ws.Range("E46").Value = "  +1.01%  "

# Row 47
$ws.Range("E47").Value = "  +0.83%  "

# Row 48
$ws.Range("D48").Value = "1.942.27"
$ws.Range("E48").Value = "  -0.27%  "

# Row 49
$ws.Range("D49").Value = "'104.98"
$ws.Range("E49").Value = "  -1.78%  "

# Row 50
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.31%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0130"
$ws.Range("E51").Value = "  -5.53%  "
